$d = $word.ActiveDocument

# --- 1. Remove the "Meta description" paragraph -----------------------
# That paragraph is the one whose text starts with "Meta description".
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^\s*Meta description") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# --- 2. Locate the DALLE image-prompt paragraph (now the last one) ----
$dalleText = "DALLE, please create a cartoon-style feature image that captures the essence of the Frozen Inferno game. The image should feature a happy Maya warrior with glasses. It can include elements such as ice and fire, wild symbols, and the game's symbols such as the castle on the rock, magic potions, and skulls. The image should be fun, visually appealing, and attention-grabbing to attract potential players."

$dalleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^\s*DALLE, please create") {
        $dalleIndex = $i
        break
    }
}

if ($dalleIndex -ne -1) {
    $dallePara = $d.Paragraphs.Item($dalleIndex)

    # --- 2a. Insert a new bold title paragraph right before it --------
    $dallePara.Range.InsertParagraphBefore()

    $newPara = $d.Paragraphs.Item($dalleIndex)
    $titleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Frozen Inferno free: Review of Unique Slot Game</w:t></w:r></w:p>"
    $newPara.Range.InsertXML($titleXml)

    # --- 2b. Replace the DALLE prompt text with the meta description --
    #         text, keeping the existing italic run formatting intact.
    $dalleIndex2 = $dalleIndex + 1
    $dallePara2 = $d.Paragraphs.Item($dalleIndex2)
    $dallePara2.Range.Find.Execute(
        $dalleText,
        $true, $false, $false, $false, $false, $true, 1, $false,
        "Discover the icy and fiery realms as you play Frozen Inferno for free. This unique online slot game features bonus features and an RTP of 96.36%.",
        2)
}
